$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 993
$ws.Range("I12").Value = 990.6667
$ws.Range("J12").Value = 1000
$ws.Range("K12").Value = 990.6667
$ws.Range("L12").Value = 1000
$ws.Range("M12").Value = -820.6667
$ws.Range("N12").Value = -1340

$ws.Range("H16").Value = 29999.5
$ws.Range("I16").Value = 49999
$ws.Range("J16").Value = 10000
$ws.Range("K16").Value = 49999
$ws.Range("L16").Value = 10000
$ws.Range("M16").Value = -49769
$ws.Range("N16").Value = -10460

$ws.Range("H17").Value = 195017.64
$ws.Range("I17").Value = 1000
$ws.Range("J17").Value = 200896.97
$ws.Range("K17").Value = 3000
$ws.Range("L17").Value = 602690.91
$ws.Range("M17").Value = -2832
$ws.Range("N17").Value = -603026.91

$ws.Range("H38").Value = 654
$ws.Range("I38").Value = 93.375
$ws.Range("J38").Value = 2149
$ws.Range("K38").Value = 280.125
$ws.Range("L38").Value = 6447
$ws.Range("M38").Value = 91.875
$ws.Range("N38").Value = -7191

$ws.Range("H70").Value = 1090.2727
$ws.Range("I70").Value = 566.3333
$ws.Range("J70").Value = 1286.75
$ws.Range("K70").Value = 1698.9999
$ws.Range("L70").Value = 3860.25
$ws.Range("M70").Value = -1428.9999
$ws.Range("N70").Value = -4400.25

$ws.Range("H73").Value = 1090.2727
$ws.Range("I73").Value = 566.3333
$ws.Range("J73").Value = 1286.75
$ws.Range("K73").Value = 1698.9999
$ws.Range("L73").Value = 3860.25
$ws.Range("M73").Value = -762.9999
$ws.Range("N73").Value = -5732.25

$ws.Range("H98").Value = 3517.8125
$ws.Range("I98").Value = 3050.5
$ws.Range("J98").Value = 3881.2778
$ws.Range("K98").Value = 3050.5
$ws.Range("L98").Value = 3881.2778
$ws.Range("M98").Value = -1552.5
$ws.Range("N98").Value = -6877.2778

$ws.Range("H122").Value = 3517.8125
$ws.Range("I122").Value = 3050.5
$ws.Range("J122").Value = 3881.2778
$ws.Range("K122").Value = 9151.5
$ws.Range("L122").Value = 11643.8334
$ws.Range("M122").Value = -6701.5
$ws.Range("N122").Value = -16543.8334

$ws.Range("H132").Value = 60231.766
$ws.Range("I132").Value = 63808.75
$ws.Range("J132").Value = 3000
$ws.Range("K132").Value = 191426.25
$ws.Range("L132").Value = 9000
$ws.Range("M132").Value = -188896.25
$ws.Range("N132").Value = -14060

$ws.Range("H137").Value = 1421.8889
$ws.Range("I137").Value = 1274.625
$ws.Range("J137").Value = 2600
$ws.Range("K137").Value = 3823.875
$ws.Range("L137").Value = 7800
$ws.Range("M137").Value = -1273.875
$ws.Range("N137").Value = -12900

$ws.Range("H138").Value = 4485.0884
$ws.Range("I138").Value = 8760.706
$ws.Range("J138").Value = 3059.8823
$ws.Range("K138").Value = 26282.118
$ws.Range("L138").Value = 9179.6469
$ws.Range("M138").Value = -21142.118
$ws.Range("N138").Value = -19459.6469

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H25").Value = 27266.666
$ws.Range("I25").Value = 1800
$ws.Range("J25").Value = 40000
$ws.Range("K25").Value = 1800
$ws.Range("L25").Value = 40000
$ws.Range("M25").Value = -1398
$ws.Range("N25").Value = -40804

$ws.Range("H32").Value = 182498.55
$ws.Range("I32").Value = 218382.78
$ws.Range("J32").Value = 17431.1
$ws.Range("K32").Value = 218382.78
$ws.Range("L32").Value = 17431.1
$ws.Range("M32").Value = -218095.78
$ws.Range("N32").Value = -18005.1

$ws.Range("H110").Value = 1726.2142
$ws.Range("I110").Value = 1817.4
$ws.Range("J110").Value = 1498.25
$ws.Range("K110").Value = 1817.4
$ws.Range("L110").Value = 1498.25
$ws.Range("M110").Value = 227.5999999999999
$ws.Range("N110").Value = -5588.25

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H76").Value = 12500
$ws.Range("I76").Value = 0
$ws.Range("J76").Value = 12500
$ws.Range("K76").Value = 0
$ws.Range("L76").Value = 12500
$ws.Range("N76").Value = -13130

$ws.Range("H79").Value = 12500
$ws.Range("I79").Value = 0
$ws.Range("J79").Value = 12500
$ws.Range("K79").Value = 0
$ws.Range("L79").Value = 12500
$ws.Range("N79").Value = -14684

$ws.Range("H80").Value = 76923270
$ws.Range("I80").Value = 166666820
$ws.Range("J80").Value = 233.14285
$ws.Range("K80").Value = 166666820
$ws.Range("L80").Value = 233.14285
$ws.Range("M80").Value = -166665822
$ws.Range("N80").Value = -2229.14285

$ws.Range("H83").Value = 76923270
$ws.Range("I83").Value = 166666820
$ws.Range("J83").Value = 233.14285
$ws.Range("K83").Value = 833334100
$ws.Range("L83").Value = 1165.71425
$ws.Range("M83").Value = -833329108
$ws.Range("N83").Value = -11149.71425

$ws.Range("H94").Value = 2704
$ws.Range("I94").Value = 1842.8334
$ws.Range("J94").Value = 4426.3335
$ws.Range("K94").Value = 1842.8334
$ws.Range("L94").Value = 4426.3335
$ws.Range("M94").Value = -1391.8334
$ws.Range("N94").Value = -5328.3335

$ws.Range("H99").Value = 9621.6875
$ws.Range("I99").Value = 12820.728
$ws.Range("J99").Value = 2583.8
$ws.Range("K99").Value = 12820.728
$ws.Range("L99").Value = 2583.8
$ws.Range("M99").Value = -11322.728
$ws.Range("N99").Value = -5579.8

$ws.Range("H134").Value = 3045
$ws.Range("I134").Value = 2910.8462
$ws.Range("J134").Value = 3219.4
$ws.Range("K134").Value = 8732.5386
$ws.Range("L134").Value = 9658.200000000001
$ws.Range("M134").Value = -6197.5386
$ws.Range("N134").Value = -14728.2

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2624.4792
$ws.Range("I31").Value = 3317.1538
$ws.Range("J31").Value = 2367.2
$ws.Range("K31").Value = 3317.1538
$ws.Range("L31").Value = 2367.2
$ws.Range("M31").Value = -3022.1538
$ws.Range("N31").Value = -2957.2

$ws.Range("H34").Value = 2624.4792
$ws.Range("I34").Value = 3317.1538
$ws.Range("J34").Value = 2367.2
$ws.Range("K34").Value = 3317.1538
$ws.Range("L34").Value = 2367.2
$ws.Range("M34").Value = -3115.1538
$ws.Range("N34").Value = -2771.2

$ws.Range("H82").Value = 0
$ws.Range("I82").Value = 0
$ws.Range("J82").Value = 0
$ws.Range("K82").Value = 0
$ws.Range("L82").Value = 0
$ws.Range("N82").ClearContents()

$ws.Range("H85").Value = 0
$ws.Range("I85").Value = 0
$ws.Range("J85").Value = 0
$ws.Range("K85").Value = 0
$ws.Range("L85").Value = 0
$ws.Range("N85").ClearContents()

$ws.Range("H138").Value = 95553
$ws.Range("I138").Value = 80000
$ws.Range("J138").Value = 99441.25
$ws.Range("K138").Value = 80000
$ws.Range("L138").Value = 99441.25
$ws.Range("M138").Value = -74860
$ws.Range("N138").Value = -109721.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H9").Value = 8337082.5
$ws.Range("I9").Value = 0
$ws.Range("J9").Value = 8337082.5
$ws.Range("K9").Value = 0
$ws.Range("L9").Value = 25011247.5
$ws.Range("M9").ClearContents()
$ws.Range("N9").Value = -25011695.5

$ws.Range("H49").Value = 4665
$ws.Range("I49").Value = 5000
$ws.Range("J49").Value = 4497.5
$ws.Range("K49").Value = 15000
$ws.Range("L49").Value = 13492.5
$ws.Range("M49").Value = -14844
$ws.Range("N49").Value = -13804.5

$ws.Range("H50").Value = 1934.8889
$ws.Range("I50").Value = 194.75
$ws.Range("J50").Value = 3327
$ws.Range("K50").Value = 584.25
$ws.Range("L50").Value = 9981
$ws.Range("M50").Value = -103.25
$ws.Range("N50").Value = -10943

$ws.Range("H53").Value = 1934.8889
$ws.Range("I53").Value = 194.75
$ws.Range("J53").Value = 3327
$ws.Range("K53").Value = 584.25
$ws.Range("L53").Value = 9981
$ws.Range("M53").Value = -103.25
$ws.Range("N53").Value = -10943

$ws.Range("H96").Value = 7366120
$ws.Range("I96").Value = 14713132
$ws.Range("J96").Value = 19108.5
$ws.Range("K96").Value = 44139396
$ws.Range("L96").Value = 57325.5
$ws.Range("M96").Value = -44137337
$ws.Range("N96").Value = -61443.5

$ws.Range("H107").Value = 826.8333
$ws.Range("I107").Value = 457
$ws.Range("J107").Value = 1196.6666
$ws.Range("K107").Value = 1371
$ws.Range("L107").Value = 3589.9998
$ws.Range("M107").Value = 549
$ws.Range("N107").Value = -7429.9998

$ws.Range("H131").Value = 5408580.5
$ws.Range("I131").Value = 9092117
$ws.Range("J131").Value = 146385.58
$ws.Range("K131").Value = 27276351
$ws.Range("L131").Value = 439156.74
$ws.Range("M131").Value = -27271311
$ws.Range("N131").Value = -449236.74

$ws.Range("H132").Value = 3631.7727
$ws.Range("I132").Value = 2992.7856
$ws.Range("J132").Value = 4750
$ws.Range("K132").Value = 26935.0704
$ws.Range("L132").Value = 42750
$ws.Range("M132").Value = -24405.0704
$ws.Range("N132").Value = -47810

$ws.Range("H140").Value = 2455.8823
$ws.Range("I140").Value = 2075

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H12").Value = 9999
$ws.Range("I12").Value = 9999
$ws.Range("J12").Value = 0
$ws.Range("K12").Value = 9999
$ws.Range("L12").Value = 0
$ws.Range("M12").Value = -9859

$ws.Range("H132").Value = 1677.9445
$ws.Range("I132").Value = 1806.6875
$ws.Range("J132").Value = 1574.95
$ws.Range("K132").Value = 5420.0625
$ws.Range("L132").Value = 4724.85
$ws.Range("M132").Value = -2890.0625
$ws.Range("N132").Value = -9784.85

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H3").Value = 5000
$ws.Range("I3").Value = 2000
$ws.Range("J3").Value = 6500
$ws.Range("K3").Value = 2000
$ws.Range("L3").Value = 6500
$ws.Range("M3").Value = -1888
$ws.Range("N3").Value = -6724

$ws.Range("H15").Value = 5000
$ws.Range("I15").Value = 2000
$ws.Range("J15").Value = 6500
$ws.Range("K15").Value = 2000
$ws.Range("L15").Value = 6500
$ws.Range("M15").Value = -1830
$ws.Range("N15").Value = -6840

$ws.Range("H16").Value = 1411.85
$ws.Range("I16").Value = 1124.3334
$ws.Range("J16").Value = 3999.5
$ws.Range("K16").Value = 1124.3334
$ws.Range("L16").Value = 3999.5
$ws.Range("M16").Value = -954.3334
$ws.Range("N16").Value = -4339.5

$ws.Range("H20").Value = 24999.5
$ws.Range("I20").Value = 25000
$ws.Range("J20").Value = 24999
$ws.Range("K20").Value = 25000
$ws.Range("L20").Value = 24999
$ws.Range("M20").Value = -24774
$ws.Range("N20").Value = -25451

$ws.Range("H132").Value = 4373.5
$ws.Range("I132").Value = 3500
$ws.Range("J132").Value = 4664.6665
$ws.Range("K132").Value = 10500
$ws.Range("L132").Value = 13993.9995
$ws.Range("M132").Value = -7970
$ws.Range("N132").Value = -19053.9995

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H7").Value = 9999
$ws.Range("I7").Value = 9999
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 9999
$ws.Range("L7").Value = 0
$ws.Range("M7").Value = -9886

$ws.Range("H81").Value = 85708.5
$ws.Range("I81").Value = 3107.7144
$ws.Range("J81").Value = 201349.6
$ws.Range("K81").Value = 6215.4288
$ws.Range("L81").Value = 402699.2
$ws.Range("M81").Value = -5154.4288
$ws.Range("N81").Value = -404821.2

$ws.Range("H84").Value = 85708.5
$ws.Range("I84").Value = 3107.7144
$ws.Range("J84").Value = 201349.6
$ws.Range("K84").Value = 31077.144
$ws.Range("L84").Value = 2013496
$ws.Range("M84").Value = -25773.144
$ws.Range("N84").Value = -2024104

$ws.Range("H123").Value = 101247.25
$ws.Range("I123").Value = 0
$ws.Range("J123").Value = 101247.25
$ws.Range("K123").Value = 0
$ws.Range("L123").Value = 101247.25
$ws.Range("N123").Value = -111047.25

$ws.Range("H126").Value = 2612.3572
$ws.Range("I126").Value = 2276.5557
$ws.Range("J126").Value = 3216.8
$ws.Range("K126").Value = 6829.6671
$ws.Range("L126").Value = 9650.400000000001
$ws.Range("M126").Value = -4359.6671
